$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(1, 1).Value = "John"
$ws.Cells.Item(1, 2).Value = 25
$ws.Cells.Item(1, 3).Value = "New York"
$ws.Cells.Item(1, 4).Value = "25th street"

$ws.Cells.Item(2, 1).Value = "Alice"
$ws.Cells.Item(2, 2).Value = 30
$ws.Cells.Item(2, 3).Value = "Los Angeles"
$ws.Cells.Item(2, 4).Value = "apartment 25, 5th floor"

$ws.Cells.Item(3, 1).Value = "Bob"
$ws.Cells.Item(3, 2).Value = 22
$ws.Cells.Item(3, 3).Value = "Chicago"
$ws.Cells.Item(3, 4).Value = "26th street"

$ws.Cells.Item(4, 1).Value = "Charlie"
$ws.Cells.Item(4, 2).Value = 28
$ws.Cells.Item(4, 3).Value = "Houston"
$ws.Cells.Item(4, 4).Value = "apartment 25, 5th floor"

$ws.Cells.Item(5, 1).Value = "David"
$ws.Cells.Item(5, 2).Value = 35
$ws.Cells.Item(5, 3).Value = "Phoenix"
$ws.Cells.Item(5, 4).Value = "27th street"

$ws.Cells.Item(6, 1).Value = "Emma"
$ws.Cells.Item(6, 2).Value = 40
$ws.Cells.Item(6, 3).Value = "Philadelphia"
$ws.Cells.Item(6, 4).Value = "apartment 25, 5th floor"

$ws.Cells.Item(7, 1).Value = "Fathi"
$ws.Cells.Item(7, 2).Value = 19
$ws.Cells.Item(7, 3).Value = "San Antonio"
$ws.Cells.Item(7, 4).Value = "28th street"

$ws.Cells.Item(8, 1).Value = "Grace"
$ws.Cells.Item(8, 2).Value = 21
$ws.Cells.Item(8, 3).Value = "San Diego"
$ws.Cells.Item(8, 4).Value = "apartment 25, 5th floor"

$ws.Cells.Item(9, 1).Value = "Henry"
$ws.Cells.Item(9, 2).Value = 45
$ws.Cells.Item(9, 3).Value = "Dallas"
$ws.Cells.Item(9, 4).Value = "29th street"

$ws.Cells.Item(10, 1).Value = "Ivy"
$ws.Cells.Item(10, 2).Value = 50
$ws.Cells.Item(10, 3).Value = "San Jose"
$ws.Cells.Item(10, 4).Value = "apartment 25, 5th floor"

$ws.Cells.Item(11, 1).Value = "Jack"
$ws.Cells.Item(11, 2).Value = 33
$ws.Cells.Item(11, 3).Value = "Austin"
$ws.Cells.Item(11, 4).Value = "30th street"

$ws.Cells.Item(12, 1).Value = "Karen"
$ws.Cells.Item(12, 2).Value = 27
$ws.Cells.Item(12, 3).Value = "Jacksonville"
$ws.Cells.Item(12, 4).Value = "apartment 25, 5th floor"

$ws.Cells.Item(13, 1).Value = "Leo"
$ws.Cells.Item(13, 2).Value = 31
$ws.Cells.Item(13, 3).Value = "Fort Worth"
$ws.Cells.Item(13, 4).Value = "31st street"

$ws.Cells.Item(14, 1).Value = "Mona"
$ws.Cells.Item(14, 2).Value = 29
$ws.Cells.Item(14, 3).Value = "Columbus"
$ws.Cells.Item(14, 4).Value = "apartment 25, 5th floor"

$ws.Cells.Item(15, 1).Value = "Nathan"
$ws.Cells.Item(15, 2).Value = 26
$ws.Cells.Item(15, 3).Value = "Charlotte"
$ws.Cells.Item(15, 4).Value = "32nd street"

$ws.Cells.Item(16, 1).Value = "Olivia"
$ws.Cells.Item(16, 2).Value = 38
$ws.Cells.Item(16, 3).Value = "Indianapolis"
$ws.Cells.Item(16, 4).Value = "apartment 25, 5th floor"

$ws.Cells.Item(17, 1).Value = "Paul"
$ws.Cells.Item(17, 2).Value = 41
$ws.Cells.Item(17, 3).Value = "Seattle"
$ws.Cells.Item(17, 4).Value = "33rd street"

$ws.Cells.Item(18, 1).Value = "Quinn"
$ws.Cells.Item(18, 2).Value = 20
$ws.Cells.Item(18, 3).Value = "Denver"
$ws.Cells.Item(18, 4).Value = "apartment 25, 5th floor"

$ws.Cells.Item(19, 1).Value = "Rachel"
$ws.Cells.Item(19, 2).Value = 32
$ws.Cells.Item(19, 3).Value = "Washington"
$ws.Cells.Item(19, 4).Value = "34th street"

$ws.Cells.Item(20, 1).Value = "Steve"
$ws.Cells.Item(20, 2).Value = 37
$ws.Cells.Item(20, 3).Value = "Boston"
$ws.Cells.Item(20, 4).Value = "apartment 25, 5th floor"
